$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: new changelog entry for "pledge read" / "Logger trace level added"
$ws.Range("A7").Value = 1.05
$ws.Range("C7").Value = "Logger trace level added"
$ws.Range("B7").Value = "pledge read"

# Row 8: new changelog entry for "arrangement read"
$ws.Range("B8").Value = "arrangement read"

# Move the active selection down to B9, matching where editing left off
$ws.Range("B9").Select()
